$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "66.094.56"
Set-TextValue $ws.Range("E2") "  +2.10%  "
Set-TextValue $ws.Range("D3") "3.236.35"
Set-TextValue $ws.Range("E3") "  +5.76%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.21%  "
Set-TextValue $ws.Range("D5") "579.81"
Set-TextValue $ws.Range("E5") "  +3.73%  "
Set-TextValue $ws.Range("D6") "150.97"
Set-TextValue $ws.Range("E6") "  +5.77%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  +0.22%  "
Set-TextValue $ws.Range("D8") "3.230.48"
Set-TextValue $ws.Range("E8") "  +5.86%  "
Set-TextValue $ws.Range("E9") "  +4.16%  "
Set-TextValue $ws.Range("D10") "7.05"
Set-TextValue $ws.Range("E10") "  +7.68%  "
Set-TextValue $ws.Range("E11") "  +4.17%  "
Set-TextValue $ws.Range("D12") "0.488"
Set-TextValue $ws.Range("E12") "  +4.58%  "
Set-TextValue $ws.Range("D13") "38.03"
Set-TextValue $ws.Range("E13") "  +4.62%  "
Set-TextValue $ws.Range("E14") "  +4.50%  "
Set-TextValue $ws.Range("D15") "3.758.14"
Set-TextValue $ws.Range("E15") "  +6.90%  "
Set-TextValue $ws.Range("D16") "66.203.30"
Set-TextValue $ws.Range("E16") "  +2.44%  "
Set-TextValue $ws.Range("D17") "539.68"
Set-TextValue $ws.Range("E17") "  +9.84%  "
Set-TextValue $ws.Range("D18") "3.239.08"
Set-TextValue $ws.Range("E18") "  +6.38%  "
Set-TextValue $ws.Range("E19") "  +2.71%  "
Set-TextValue $ws.Range("D20") "7.13"
Set-TextValue $ws.Range("E20") "  +6.21%  "
Set-TextValue $ws.Range("D21") "14.52"
Set-TextValue $ws.Range("E21") "  +5.77%  "
Set-TextValue $ws.Range("E22") "  +7.42%  "
Set-TextValue $ws.Range("D23") "7.78"
Set-TextValue $ws.Range("E23") "  +8.77%  "
Set-TextValue $ws.Range("D24") "13.48"
Set-TextValue $ws.Range("E24") "  +5.93%  "
Set-TextValue $ws.Range("D25") "80.97"
Set-TextValue $ws.Range("E25") "  +2.32%  "
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  +0.05%  "
Set-TextValue $ws.Range("D27") "9.24"
Set-TextValue $ws.Range("E27") "  +17.31%  "
Set-TextValue $ws.Range("E28") "  +7.99%  "
Set-TextValue $ws.Range("D29") "2.26"
Set-TextValue $ws.Range("E29") "  +6.53%  "
Set-TextValue $ws.Range("D30") "27.66"
Set-TextValue $ws.Range("E30") "  +5.96%  "
Set-TextValue $ws.Range("D31") "2.74"
Set-TextValue $ws.Range("E31") "  +3.40%  "
Set-TextValue $ws.Range("D32") "1.00"
Set-TextValue $ws.Range("E32") "  +0.21%  "
Set-TextValue $ws.Range("E33") "  +5.30%  "
Set-TextValue $ws.Range("D34") "562.92"
Set-TextValue $ws.Range("E34") "  +3.04%  "
Set-TextValue $ws.Range("E35") "  +6.15%  "
Set-TextValue $ws.Range("E36") "  +2.76%  "
Set-TextValue $ws.Range("D37") "0.0454"
Set-TextValue $ws.Range("E37") "  +7.93%  "
Set-TextValue $ws.Range("E38") "  +3.40%  "
Set-TextValue $ws.Range("D39") "0.0859"
Set-TextValue $ws.Range("E39") "  +6.26%  "
Set-TextValue $ws.Range("E40") "  +5.79%  "
Set-TextValue $ws.Range("D41") "3.186.86"
Set-TextValue $ws.Range("E41") "  +9.58%  "
Set-TextValue $ws.Range("D42") "2.90"
Set-TextValue $ws.Range("E42") "  +3.59%  "
Set-TextValue $ws.Range("D43") "8.55"
Set-TextValue $ws.Range("E43") "  +2.94%  "
Set-TextValue $ws.Range("D44") "0.288"
Set-TextValue $ws.Range("E44") "  +16.85%  "
Set-TextValue $ws.Range("D45") "2.32"
Set-TextValue $ws.Range("E45") "  +9.27%  "
Set-TextValue $ws.Range("D46") "26.38"
Set-TextValue $ws.Range("E46") "  +5.53%  "
Set-TextValue $ws.Range("E47") "  +0.03%  "
Set-TextValue $ws.Range("D48") "0.0₃0553"
Set-TextValue $ws.Range("E48") "  +2.00%  "
Set-TextValue $ws.Range("D49") "125.24"
Set-TextValue $ws.Range("E49") "  +3.98%  "
Set-TextValue $ws.Range("E50") "  +2.70%  "
Set-TextValue $ws.Range("D51") "2.20"
Set-TextValue $ws.Range("E51") "  +6.91%  "
